# "fix wrong bitstream order"
#
# The "Comment" column (G) for the en_bus_* rows (108-139) had several
# blank/placeholder ("?") cells because the bitstream bit order documented
# there was wrong. This fills in the corrected bit names, and marks the
# older "Function Name" column (F) values for that same range as
# struck-through (they are superseded by the corrected G column values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bitstream")

# ---------------------------------------------------------------------
# 1. Fill in the corrected "Comment" values in column G for rows 116-139.
# ---------------------------------------------------------------------

# Rows 116-123 were blank; fill with the en_bus_north[3..0] / en_bus_south[3..0]
# names (previously undocumented for this nibble).
$ws.Cells.Item(116, 7).Value2 = "en_bus_north[3]"
$ws.Cells.Item(117, 7).Value2 = "en_bus_north[2]"
$ws.Cells.Item(118, 7).Value2 = "en_bus_north[1]"
$ws.Cells.Item(119, 7).Value2 = "en_bus_north[0]"
$ws.Cells.Item(120, 7).Value2 = "en_bus_south[3]"
$ws.Cells.Item(121, 7).Value2 = "en_bus_south[2]"
$ws.Cells.Item(122, 7).Value2 = "en_bus_south[1]"
$ws.Cells.Item(123, 7).Value2 = "en_bus_south[0]"

# Rows 124-131 were blank; fill with the en_bus_east[3..0] / en_bus_west[3..0]
# names.
$ws.Cells.Item(124, 7).Value2 = "en_bus_east[3]"
$ws.Cells.Item(125, 7).Value2 = "en_bus_east[2]"
$ws.Cells.Item(126, 7).Value2 = "en_bus_east[1]"
$ws.Cells.Item(127, 7).Value2 = "en_bus_east[0]"
$ws.Cells.Item(128, 7).Value2 = "en_bus_west[3]"
$ws.Cells.Item(129, 7).Value2 = "en_bus_west[2]"
$ws.Cells.Item(130, 7).Value2 = "en_bus_west[1]"
$ws.Cells.Item(131, 7).Value2 = "en_bus_west[0]"

# Rows 132-139 held a placeholder "?"; replace with the real bit[5]/bit[4]
# enable names.
$ws.Cells.Item(132, 7).Value2 = "en_bus_east[5]"
$ws.Cells.Item(133, 7).Value2 = "en_bus_east[4]"
$ws.Cells.Item(134, 7).Value2 = "en_bus_west[5]"
$ws.Cells.Item(135, 7).Value2 = "en_bus_west[4]"
$ws.Cells.Item(136, 7).Value2 = "en_bus_north[5]"
$ws.Cells.Item(137, 7).Value2 = "en_bus_north[4]"
$ws.Cells.Item(138, 7).Value2 = "en_bus_south[5]"
$ws.Cells.Item(139, 7).Value2 = "en_bus_south[4]"

# ---------------------------------------------------------------------
# 2. Give the newly-populated G cells the same (Consolas) formatting
#    already used by the other Comment cells in this table, by copying
#    the format from an already-styled comment cell (G112).
# ---------------------------------------------------------------------
$fmtSrc = $ws.Cells.Item(112, 7)
$fmtSrc.Copy() | Out-Null
$ws.Range("G116:G139").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Strike through the now-superseded "Function Name" values in column F
#    for the whole en_bus_* block (rows 108-139), to show they've been
#    replaced by the corrected Comment column.
# ---------------------------------------------------------------------
for ($r = 108; $r -le 139; $r++) {
    $ws.Cells.Item($r, 6).Font.Strikethrough = $true
}

# ---------------------------------------------------------------------
# 4. Restore the cursor/selection to where the author left it.
# ---------------------------------------------------------------------
$ws.Range("G109").Select() | Out-Null
